$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep their original Text storage type (the source data
# stores Coin/Link/Price/Volume columns as text, e.g. "34.20" or "0.0230"); without
# forcing a Text number format first, Excel auto-detects numeric-looking strings and
# would silently drop significant trailing zeros (e.g. "34.20" -> 34.2).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.952.49'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.453.62'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.33%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '490.38'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.71'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +10.46%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.80%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.457.99'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.68%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.60%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.38%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.83%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.66%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.880.75'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.54%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.151.61'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.01'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.59%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000138'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.55%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.468.42'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.00%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.94%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '324.98'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.47%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.82%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.32%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.16'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.91%  '

$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.01'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.00%  '

$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.407'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.84%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.163'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.566.13'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.55'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.44%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0815'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +7.32%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.05%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.94'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.11%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.53'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.78%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.22'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.26'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.77%  '

$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.15'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.36%  '

$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.887'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.97%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.76'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.41'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +10.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.20'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.51%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.20%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0559'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.20%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.26%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.607'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0965'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +8.11%  '

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.84'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.87%  '

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '264.25'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.15%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0230'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.04%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.24'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.96%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.82'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.43%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.75'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +27.32%  '
